# Auto-generated edit script applying the diff changes to before.xlsx
# Update "想去人数" (F column) counters across all 4 sheets, and
# replace row 2 content on the "演出" sheet with a new event.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('展览')
$ws.Range('F3').Value = 3617
$ws.Range('F5').Value = 8420
$ws.Range('F7').Value = 145
$ws.Range('F8').Value = 2313
$ws.Range('F10').Value = 127
$ws.Range('F11').Value = 7690
$ws.Range('F13').Value = 5043
$ws.Range('F17').Value = 5408
$ws.Range('F20').Value = 148
$ws.Range('F22').Value = 976
$ws.Range('F23').Value = 1514
$ws.Range('F24').Value = 2127
$ws.Range('F25').Value = 30
$ws.Range('F26').Value = 213
$ws.Range('F27').Value = 269
$ws.Range('F28').Value = 1107
$ws.Range('F30').Value = 762
$ws.Range('F33').Value = 1311
$ws.Range('F34').Value = 492
$ws.Range('F35').Value = 25
$ws.Range('F37').Value = 237
$ws.Range('F38').Value = 29
$ws.Range('F39').Value = 69
$ws.Range('F41').Value = 2514

$ws = $wb.Worksheets.Item('演出')
$ws.Range('F3').Value = 7832
$ws.Range('F31').Value = 85
$ws.Range('F41').Value = 158

$ws = $wb.Worksheets.Item('本地生活')
$ws.Range('F4').Value = 2424
$ws.Range('F7').Value = 699
$ws.Range('F9').Value = 9497
$ws.Range('F10').Value = 1827
$ws.Range('F15').Value = 315
$ws.Range('F16').Value = 2626
$ws.Range('F17').Value = 302
$ws.Range('F18').Value = 114
$ws.Range('F19').Value = 579

$ws = $wb.Worksheets.Item('全部类型')
$ws.Range('F3').Value = 3617
$ws.Range('F4').Value = 2424
$ws.Range('F5').Value = 1827
$ws.Range('F7').Value = 315
$ws.Range('F8').Value = 2626
$ws.Range('F9').Value = 302
$ws.Range('F10').Value = 127
$ws.Range('F11').Value = 7690
$ws.Range('F15').Value = 148
$ws.Range('F16').Value = 114
$ws.Range('F17').Value = 976
$ws.Range('F18').Value = 1514
$ws.Range('F19').Value = 2127
$ws.Range('F21').Value = 579
$ws.Range('F22').Value = 579
$ws.Range('F26').Value = 269
$ws.Range('F27').Value = 762
$ws.Range('F31').Value = 1311
$ws.Range('F35').Value = 492
$ws.Range('F39').Value = 237
$ws.Range('F46').Value = 2514

# 演出 sheet row 2: swap out old event for the new "炒饭·二次元律动1st Live" event
$ws = $wb.Worksheets.Item('演出')
# Force B2 to stay a text value (not auto-parsed into a date serial) just like
# the original "2024-10-04" inline string, then restore General display format.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = '2024-10-05'
$ws.Range("B2").NumberFormat = "General"
$ws.Range("C2").Value = '上海·“炒饭·二次元律动1st Live”'
$ws.Range("D2").Value = '衡山路八号水塔广场 JZ Club 爵士上海俱乐部'
$ws.Range("E2").Value = '2024.10.05 16:00-10.05 18:00'
$ws.Range("F2").Value = 60
$ws.Range("G2").Value = '已售罄'
$ws.Range("H2").Value = 'https://show.bilibili.com/platform/detail.html?id=92183'
$ws.Range("I2").Value = '//i2.hdslb.com/bfs/openplatform/202409/PeGFMPZC1725868905755.jpeg'
